$p = $ppt.ActivePresentation

# --- Slide 1: title shape ---
$s1 = $p.Slides.Item(1)
$title = $s1.Shapes.Item(1)
$title.Left = 609600 / 12700
$title.Top = 533400 / 12700
$title.Width = 7772400 / 12700
$title.Height = 1470025 / 12700
# The final text is identical (once concatenated) to the existing two-run
# text, so set an intermediate placeholder first to force the two runs to
# collapse into a single run, then set the final text.
$title.TextFrame.TextRange.Text = "placeholder_tmp"
$title.TextFrame.TextRange.Text = "How to make an app for Zendesk"
